# "cleaned up code & added comments"
#
# Slide 3 ("Main Functionality") content placeholder: the first bullet
# "Use of Arduino Mega to act as main controller and data acquisition
# system" is reworded to "Use of (2) Arduino Mega microcontroller
# boards" and split into a top-level bullet with two new sub-bullets
# ("main controller" / "data acquisition system"). The placeholder also
# picks up PowerPoint's normal autofit shrink (text-to-fit-shape) since
# it now needs to fit more lines.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Re-word the first bullet. Assigning the final string directly would
# make the host diff against the old text and keep the shared "Use of "
# prefix as a separate leftover run, so first stomp it with an
# unrelated placeholder (no shared prefix) and then assign the real
# text, which collapses back down to a single run.
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "X"
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "Use of (2) Arduino Mega microcontroller boards"

# Insert the two new sub-bullets ahead of the "Load cell..." bullet.
$para2 = $tr.Paragraphs(2, 1)
$para2.InsertBefore("main controller`rdata acquisition system`r")

# Demote the two freshly-inserted paragraphs to the second outline level.
$newPara2 = $tr.Paragraphs(2, 1)
$newPara2.IndentLevel = 2
$newPara3 = $tr.Paragraphs(3, 1)
$newPara3.IndentLevel = 2

# The placeholder now has more lines of text than before, so turn on
# PowerPoint's "shrink text on overflow" autofit for it.
$tf.AutoSize = 2
